$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 90
$ws.Range("B3").Value = 50
$ws.Range("C3").Value = "get_first_avail"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0.4
$ws.Range("F3").Value = 0.25
$ws.Range("G3").Value = 0.15
$ws.Range("H3").Value = 0.1
$ws.Range("I3").Value = 0.1

$ws.Range("A4").Select()
